$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 463 (currently "A_SEA_B") to add the new
# "K_RISIKO" entries, shifting existing rows 463-535 down to 467-539.
$insertRange = $ws.Range("A463:D466")
$insertRange.Insert()

# Copy formatting (style) from the row that now follows the inserted block
# (row 467, formerly row 463) down onto the 4 new rows.
$formatSource = $ws.Range("A467:D467")
$formatSource.Copy()
$ws.Range("A463:D466").PasteSpecial(-4122)

# Row 463: A_RISIKO_0
$ws.Cells.Item(463, 1).Value2 = "A_RISIKO_0"
$ws.Cells.Item(463, 2).Value2 = "K_RISIKO"
$ws.Cells.Item(463, 3).Value2 = "Keine Risikolage"
$ws.Cells.Item(463, 4).Value2 = "XXXKeine Risikolage"

# Row 464: A_RISIKO_1
$ws.Cells.Item(464, 1).Value2 = "A_RISIKO_1"
$ws.Cells.Item(464, 2).Value2 = "K_RISIKO"
$ws.Cells.Item(464, 3).Value2 = "Mindestens eine Risikolage"
$ws.Cells.Item(464, 4).Value2 = "XXXMindestens eine Risikolage"

# Row 465: A_RISIKO_ALL
$ws.Cells.Item(465, 1).Value2 = "A_RISIKO_ALL"
$ws.Cells.Item(465, 2).Value2 = "K_RISIKO"
$ws.Cells.Item(465, 3).Value2 = "Insgesamt"
$ws.Cells.Item(465, 4).Value2 = "Total"

# Row 466: A_RISIKO_DIFF
$ws.Cells.Item(466, 1).Value2 = "A_RISIKO_DIFF"
$ws.Cells.Item(466, 2).Value2 = "K_RISIKO"
$ws.Cells.Item(466, 3).Value2 = "Differenz zwischen `"Keine Risikolage`" und `"Mind. eine Risikolage`""
$ws.Cells.Item(466, 4).Value2 = "XXXDifferenz zwischen `"Keine Risikolage`" und `"Mind. eine Risikolage`""
